# Update countries & provincias Spain
# Refresh of the live COVID country stats: the timestamp advances from
# 03:52 to 04:22 and several countries receive new totals. A handful of
# countries with very close case counts swap rank order once the table is
# re-sorted descending by "Casos totales", which is why both the country
# name and the numbers change together for those rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header timestamp
$ws.Range("A1").Value = "Datos actualizados a 16 de Abril de 2020 a las 04:22"

# Noruega keeps its row but gets fresh totals
$ws.Range("B32").Value = 6798
$ws.Range("C32").Value = 1
$ws.Range("D32").Value = 32
$ws.Range("E32").Value = 6616
$ws.Range("F32").Value = 64
$ws.Range("G32").Value = 0
$ws.Range("H32").Value = 150

# Crucero keeps its row but gets fresh totals
$ws.Range("B87").Value = 712
$ws.Range("C87").Value = 0
$ws.Range("D87").Value = 644
$ws.Range("E87").Value = 56
$ws.Range("F87").Value = 7
$ws.Range("G87").Value = 0
$ws.Range("H87").Value = 12

# Rows 99-105: Republica de Yibuti / Bolivia swap rank, everything refreshed
$ws.Range("A99").Value = "Bolivia"
$ws.Range("B99").Value = 441
$ws.Range("C99").Value = 44
$ws.Range("D99").Value = 14
$ws.Range("E99").Value = 398
$ws.Range("F99").Value = 3
$ws.Range("G99").Value = 1
$ws.Range("H99").Value = 29

$ws.Range("A100").Value = "Republica de Yibuti"
$ws.Range("B100").Value = 435
$ws.Range("C100").Value = 0
$ws.Range("D100").Value = 71
$ws.Range("E100").Value = 362
$ws.Range("F100").Value = 0
$ws.Range("G100").Value = 0
$ws.Range("H100").Value = 2

$ws.Range("A101").Value = "Honduras"
$ws.Range("B101").Value = 426
$ws.Range("C101").Value = 7
$ws.Range("D101").Value = 9
$ws.Range("E101").Value = 382
$ws.Range("F101").Value = 10
$ws.Range("G101").Value = 4
$ws.Range("H101").Value = 35

$ws.Range("A102").Value = "Nigeria"
$ws.Range("B102").Value = 407
$ws.Range("C102").Value = 0
$ws.Range("D102").Value = 128
$ws.Range("E102").Value = 267
$ws.Range("F102").Value = 2
$ws.Range("G102").Value = 0
$ws.Range("H102").Value = 12

$ws.Range("A103").Value = "Guinea"
$ws.Range("B103").Value = 404
$ws.Range("C103").Value = 0
$ws.Range("D103").Value = 31
$ws.Range("E103").Value = 372
$ws.Range("F103").Value = 0
$ws.Range("G103").Value = 0
$ws.Range("H103").Value = 1

$ws.Range("A104").Value = "Jordania"
$ws.Range("B104").Value = 401
$ws.Range("C104").Value = 0
$ws.Range("D104").Value = 250
$ws.Range("E104").Value = 144
$ws.Range("F104").Value = 5
$ws.Range("G104").Value = 0
$ws.Range("H104").Value = 7

$ws.Range("A105").Value = "Malta"
$ws.Range("B105").Value = 399
$ws.Range("C105").Value = 0
$ws.Range("D105").Value = 44
$ws.Range("E105").Value = 352
$ws.Range("F105").Value = 4
$ws.Range("G105").Value = 0
$ws.Range("H105").Value = 3

# Rows 121-122: Islas Feroe / Guatemala swap rank
$ws.Range("A121").Value = "Guatemala"
$ws.Range("B121").Value = 196
$ws.Range("C121").Value = 16
$ws.Range("D121").Value = 19
$ws.Range("E121").Value = 172
$ws.Range("F121").Value = 3
$ws.Range("G121").Value = 0
$ws.Range("H121").Value = 5

$ws.Range("A122").Value = "Islas Feroe"
$ws.Range("B122").Value = 184
$ws.Range("C122").Value = 0
$ws.Range("D122").Value = 166
$ws.Range("E122").Value = 18
$ws.Range("F122").Value = 0
$ws.Range("G122").Value = 0
$ws.Range("H122").Value = 0

# Rows 154-155: Bahamas / San Martin (Parte Holandesa) swap rank
$ws.Range("A154").Value = "San Martin (Parte Holandesa)"
$ws.Range("B154").Value = 53
$ws.Range("C154").Value = 0
$ws.Range("D154").Value = 5
$ws.Range("E154").Value = 39
$ws.Range("F154").Value = 2
$ws.Range("G154").Value = 0
$ws.Range("H154").Value = 9

$ws.Range("A155").Value = "Bahamas"
$ws.Range("B155").Value = 53
$ws.Range("C155").Value = 0
$ws.Range("D155").Value = 6
$ws.Range("E155").Value = 39
$ws.Range("F155").Value = 1
$ws.Range("G155").Value = 0
$ws.Range("H155").Value = 8

# Rows 195-196: Montserrat / Islas Malvinas swap rank
$ws.Range("A195").Value = "Islas Malvinas"
$ws.Range("B195").Value = 11
$ws.Range("C195").Value = 0
$ws.Range("D195").Value = 1
$ws.Range("E195").Value = 10
$ws.Range("F195").Value = 0
$ws.Range("G195").Value = 0
$ws.Range("H195").Value = 0

$ws.Range("A196").Value = "Montserrat"
$ws.Range("B196").Value = 11
$ws.Range("C196").Value = 0
$ws.Range("D196").Value = 1
$ws.Range("E196").Value = 10
$ws.Range("F196").Value = 1
$ws.Range("G196").Value = 0
$ws.Range("H196").Value = 0

# Rows 209-210: Santo Tome y Principe / Sudan del Sur swap rank (totals tie, no number change)
$ws.Range("A209").Value = "Sudan del Sur"
$ws.Range("B209").Value = 4
$ws.Range("C209").Value = 0
$ws.Range("D209").Value = 0
$ws.Range("E209").Value = 4
$ws.Range("F209").Value = 0
$ws.Range("G209").Value = 0
$ws.Range("H209").Value = 0

$ws.Range("A210").Value = "Santo Tome y Principe"
$ws.Range("B210").Value = 4
$ws.Range("C210").Value = 0
$ws.Range("D210").Value = 0
$ws.Range("E210").Value = 4
$ws.Range("F210").Value = 0
$ws.Range("G210").Value = 0
$ws.Range("H210").Value = 0

# Rows 215-216: Yemen / San Pedro y Miquelon swap rank (totals tie, no number change)
$ws.Range("A215").Value = "San Pedro y Miquelon"
$ws.Range("B215").Value = 1
$ws.Range("C215").Value = 0
$ws.Range("D215").Value = 0
$ws.Range("E215").Value = 1
$ws.Range("F215").Value = 0
$ws.Range("G215").Value = 0
$ws.Range("H215").Value = 0

$ws.Range("A216").Value = "Yemen"
$ws.Range("B216").Value = 1
$ws.Range("C216").Value = 0
$ws.Range("D216").Value = 0
$ws.Range("E216").Value = 1
$ws.Range("F216").Value = 0
$ws.Range("G216").Value = 0
$ws.Range("H216").Value = 0
